$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicate "Contact" row (old row 11), leaving a single
# Publisher/Contact pair behind which we then turn into
# Publisher/Jurisdiction rows below.
$ws.Rows.Item(11).Delete()

# Publisher row gets a value, and the old "Contact" row becomes "Jurisdiction"
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive row now has a value of true, stored as the literal text
# "true" (not a Boolean). A bare Value = "true" assignment gets
# auto-coerced to a Boolean by Excel, so compute it as a text formula on
# a scratch cell, then paste only the resulting value into B14 so the
# existing cell style/formatting is left untouched.
$scratch = $ws.Range("Z100")
$scratch.Formula = "=""true"""
$scratch.Copy()
$ws.Range("B14").PasteSpecial(-4163)
$scratch.Clear()

Write-Output "done"
